$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell from "gen_id" to "genero"
$ws.Range("A1").Value = "genero"
